# Scheduled-runner refresh of market-price / profit columns (H:N) across the
# Hades_Profits leve tables. Source data comes from an external price API, so
# values are written as literals (no formulas) on each affected sheet/row.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 125
$ws.Range("I125").Value = 3000
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 27000
$ws.Range("L125").Value = 27000
$ws.Range("M125").Value = -24540

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6167.83
$ws.Range("I32").Value = 5203.043
$ws.Range("J32").Value = 18985.715
$ws.Range("K32").Value = 5203.043
$ws.Range("L32").Value = 18985.715
$ws.Range("M32").Value = -4916.043

# Row 45
$ws.Range("H45").Value = 938.7692
$ws.Range("I45").Value = 845.8182
$ws.Range("J45").Value = 1450
$ws.Range("K45").Value = 845.8182
$ws.Range("L45").Value = 1450
$ws.Range("M45").Value = -468.8182
$ws.Range("N45").Value = -2204

# Row 61
$ws.Range("H61").Value = 27082426
$ws.Range("I61").Value = 38501036
$ws.Range("J61").Value = 92980.91
$ws.Range("K61").Value = 38501036
$ws.Range("L61").Value = 92980.91
$ws.Range("M61").Value = -38500824
$ws.Range("N61").Value = -93404.91

# Row 82
$ws.Range("H82").Value = 20000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 20000
$ws.Range("N82").Value = -20722

# Row 85
$ws.Range("H85").Value = 20000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 20000
$ws.Range("N85").Value = -22496

# Row 136
$ws.Range("H136").Value = 27082426
$ws.Range("I136").Value = 38501036
$ws.Range("J136").Value = 92980.91
$ws.Range("K136").Value = 115503108
$ws.Range("L136").Value = 278942.73
$ws.Range("M136").Value = -115500558
$ws.Range("N136").Value = -284042.73

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 38712
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 38712
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 38712
$ws.Range("N2").Value = -38938

# Row 134
$ws.Range("H134").Value = 2109.4263
$ws.Range("I134").Value = 2062.7192
$ws.Range("J134").Value = 2775
$ws.Range("K134").Value = 6188.1576
$ws.Range("L134").Value = 8325
$ws.Range("M134").Value = -3653.1576
$ws.Range("N134").Value = -13395

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1309.4706
$ws.Range("I16").Value = 1113.6364
$ws.Range("J16").Value = 1668.5
$ws.Range("K16").Value = 1113.6364
$ws.Range("L16").Value = 1668.5
$ws.Range("M16").Value = -826.6364000000001
$ws.Range("N16").Value = -2242.5

# Row 48
$ws.Range("H48").Value = 3000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 3000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 3000
$ws.Range("N48").Value = -3952

# Row 99
$ws.Range("H99").Value = 1736
$ws.Range("I99").Value = 1696.6666
$ws.Range("J99").Value = 1795
$ws.Range("K99").Value = 1696.6666
$ws.Range("L99").Value = 1795
$ws.Range("M99").Value = -198.6666
$ws.Range("N99").Value = -4791

# Row 105
$ws.Range("H105").Value = 1283.5333
$ws.Range("I105").Value = 1220.1666
$ws.Range("J105").Value = 1537
$ws.Range("K105").Value = 1220.1666
$ws.Range("L105").Value = 1537
$ws.Range("M105").Value = 526.8334
$ws.Range("N105").Value = -5031

# Row 113
$ws.Range("H113").Value = 1309.4706
$ws.Range("I113").Value = 1113.6364
$ws.Range("J113").Value = 1668.5
$ws.Range("K113").Value = 1113.6364
$ws.Range("L113").Value = 1668.5
$ws.Range("M113").Value = 1056.3636
$ws.Range("N113").Value = -6008.5

# Row 115
$ws.Range("H115").Value = 25214.285
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 25214.285
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 25214.285
$ws.Range("N115").Value = -27564.285

# Row 126
$ws.Range("H126").Value = 1736
$ws.Range("I126").Value = 1696.6666
$ws.Range("J126").Value = 1795
$ws.Range("K126").Value = 5089.9998
$ws.Range("L126").Value = 5385
$ws.Range("M126").Value = -2619.9998
$ws.Range("N126").Value = -10325

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 2943.5
$ws.Range("I75").Value = 2943.5
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 8830.5
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -7832.5
$ws.Range("N75").ClearContents()

# Row 78
$ws.Range("H78").Value = 2943.5
$ws.Range("I78").Value = 2943.5
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 26491.5
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -21499.5
$ws.Range("N78").ClearContents()

# Row 108
$ws.Range("H108").Value = 404.8
$ws.Range("I108").Value = 404.8
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 1214.4
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 1665.6

# Row 141
$ws.Range("H141").Value = 9377.85
$ws.Range("I141").Value = 3229.75
$ws.Range("J141").Value = 18600
$ws.Range("K141").Value = 9689.25
$ws.Range("L141").Value = 55800
$ws.Range("M141").Value = -4509.25
$ws.Range("N141").Value = -66160

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 9750
$ws.Range("I52").Value = 9500
$ws.Range("J52").Value = 10000
$ws.Range("K52").Value = 9500
$ws.Range("L52").Value = 10000
$ws.Range("M52").Value = -9241
$ws.Range("N52").Value = -10518

# Row 113
$ws.Range("H113").Value = 2237.6667
$ws.Range("I113").Value = 1612.25
$ws.Range("J113").Value = 2738
$ws.Range("K113").Value = 1612.25
$ws.Range("L113").Value = 2738
$ws.Range("M113").Value = 557.75
$ws.Range("N113").Value = -7078

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 13893275
$ws.Range("I7").Value = 16670463
$ws.Range("J7").Value = 7335
$ws.Range("K7").Value = 16670463
$ws.Range("L7").Value = 7335
$ws.Range("M7").Value = -16670351
$ws.Range("N7").Value = -7559

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# Row 117
$ws.Range("H117").Value = 49696
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 49696
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 49696
$ws.Range("N117").Value = -58874

# Row 126
$ws.Range("H126").Value = 13893275
$ws.Range("I126").Value = 16670463
$ws.Range("J126").Value = 7335
$ws.Range("K126").Value = 50011389
$ws.Range("L126").Value = 22005
$ws.Range("M126").Value = -50008919
$ws.Range("N126").Value = -26945

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 3000
$ws.Range("I96").Value = 4000
$ws.Range("J96").Value = 2500
$ws.Range("K96").Value = 4000
$ws.Range("L96").Value = 2500
$ws.Range("M96").Value = -2627
$ws.Range("N96").Value = -5246

# Row 100
$ws.Range("H100").Value = 92027.55
$ws.Range("I100").Value = 72485.86
$ws.Range("J100").Value = 126225.5
$ws.Range("K100").Value = 144971.72
$ws.Range("L100").Value = 252451
$ws.Range("M100").Value = -144430.72
$ws.Range("N100").Value = -253533

# Row 126
$ws.Range("H126").Value = 1354.7778
$ws.Range("I126").Value = 1211.625
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 3634.875
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -1164.875
$ws.Range("N126").Value = -12440

# Row 128
$ws.Range("H128").Value = 52711.11
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 52711.11
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 52711.11
$ws.Range("N128").Value = -62671.11
